$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "33.995.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.777.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.547"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.32%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("E10").Value = "  +4.16%  "
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.031.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.781.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.621"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.973.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.61%  "
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0516"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("E33").Value = "  -3.76%  "
$ws.Range("E34").Value = "  -4.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.395.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.628"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.931"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.53%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  -5.60%  "
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0491"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("E45").Value = "  -3.17%  "
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.927.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0120"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.79%  "

Write-Host "Applied 88 cell updates"
